# Generate Report for Handoff
# - Sets Priority ("ht") for a batch of files that were just handed off
#   on both the "zh-cn" and "de-de" report sheets.
# - Refreshes the associated handoff timestamps (Overview "Latest HO Xliff
#   Generate Date", zh-cn "Latest Handoff Datetime", de-de "Latest Handoff
#   Datetime") to reflect the new handoff run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Rows (on the zh-cn / de-de sheets) for the files that were handed off
# in this run.
$rows = @(8, 9, 10, 11, 12, 14)

foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}

# Update the handoff timestamps for the same set of rows.
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-09-05 12:25:38"
    $dede.Range("H$r").Value = "2016-09-05 12:25:38"
    $zhcn.Range("H$r").Value = "2016-09-05 12:25:32"
}
